$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Nuevos documentos agregados al corpus (filas 15 y 16)
$ws.Range("A15").Value = 13
$ws.Range("B15").Value = "P_IFT_070218_83_AccUPR"
$ws.Range("C15").Value = "Victoria"

$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "Acta_43aOrd_011117_Acc"
$ws.Range("C16").Value = "Victoria"

# La columna B se ensanchó levemente al ajustar el contenido
$ws.Columns.Item(2).ColumnWidth = 20.1665

# Estado final de selección tras capturar los nuevos datos
$ws.Range("C17").Select()
